$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 47119
$ws.Range("E2").Value = 6591
$ws.Range("F2").Value = 6591
$ws.Range("G2").Value = 6818
$ws.Range("H2").Value = 4974
$ws.Range("I2").Value = 2224
$ws.Range("J2").Value = 2751
$ws.Range("K2").Value = 54997
$ws.Range("L2").Value = 11764
$ws.Range("M2").Value = 43234
$ws.Range("N2").Value = 23178
$ws.Range("O2").Value = 20055
$ws.Range("P2").Value = 445
$ws.Range("Q2").Value = 6914
$ws.Range("R2").Value = -5945
$ws.Range("S2").Value = -437
$ws.Range("T2").Value = 3326
$ws.Range("U2").Value = 3588
$ws.Range("V2").Value = 1683
$ws.Range("W2").Value = 13.99
$ws.Range("X2").Value = 10.56
$ws.Range("Y2").Value = 10.01
$ws.Range("Z2").Value = 9.550000000000001
$ws.Range("AA2").Value = 27.21
$ws.Range("AB2").Value = 5453.82
$ws.Range("AC2").Value = 2501
$ws.Range("AD2").Value = 39.81
$ws.Range("AE2").Value = 28017
$ws.Range("AF2").Value = 3.55
$ws.Range("AG2").Value = 325
$ws.Range("AH2").Value = 0.33
$ws.Range("AI2").Value = 12.19
$ws.Range("AJ2").Value = 79790980
$ws.Range("D3").Value = 56612
$ws.Range("E3").Value = 9136
$ws.Range("F3").Value = 9136
$ws.Range("G3").Value = 9216
$ws.Range("H3").Value = 6739
$ws.Range("I3").Value = 2604
$ws.Range("J3").Value = 4135
$ws.Range("K3").Value = 61480
$ws.Range("L3").Value = 12602
$ws.Range("M3").Value = 48878
$ws.Range("N3").Value = 25327
$ws.Range("O3").Value = 23551
$ws.Range("P3").Value = 445
$ws.Range("Q3").Value = 7365
$ws.Range("R3").Value = -833
$ws.Range("S3").Value = -707
$ws.Range("T3").Value = 2778
$ws.Range("U3").Value = 4587
$ws.Range("V3").Value = 1731
$ws.Range("W3").Value = 16.14
$ws.Range("X3").Value = 11.91
$ws.Range("Y3").Value = 10.74
$ws.Range("Z3").Value = 11.57
$ws.Range("AA3").Value = 25.78
$ws.Range("AB3").Value = 5947.41
$ws.Range("AC3").Value = 2929
$ws.Range("AD3").Value = 50.4
$ws.Range("AE3").Value = 30408
$ws.Range("AF3").Value = 4.86
$ws.Range("AG3").Value = 390
$ws.Range("AH3").Value = 0.26
$ws.Range("AI3").Value = 12.48
$ws.Range("AJ3").Value = 79790980
$ws.Range("D4").Value = 66976
$ws.Range("E4").Value = 10828
$ws.Range("F4").Value = 10828
$ws.Range("G4").Value = 10883
$ws.Range("H4").Value = 8115
$ws.Range("I4").Value = 3424
$ws.Range("J4").Value = 4691
$ws.Range("K4").Value = 70884
$ws.Range("L4").Value = 15215
$ws.Range("M4").Value = 55669
$ws.Range("N4").Value = 28332
$ws.Range("O4").Value = 27336
$ws.Range("P4").Value = 445
$ws.Range("Q4").Value = 8767
$ws.Range("R4").Value = -8141
$ws.Range("S4").Value = -989
$ws.Range("T4").Value = 6023
$ws.Range("U4").Value = 2744
$ws.Range("V4").Value = 2087
$ws.Range("W4").Value = 16.17
$ws.Range("X4").Value = 12.12
$ws.Range("Y4").Value = 12.76
$ws.Range("Z4").Value = 12.26
$ws.Range("AA4").Value = 27.33
$ws.Range("AB4").Value = 6628.95
$ws.Range("AC4").Value = 3851
$ws.Range("AD4").Value = 34.45
$ws.Range("AE4").Value = 34017
$ws.Range("AF4").Value = 3.9
$ws.Range("AG4").Value = 460
$ws.Range("AH4").Value = 0.35
$ws.Range("AI4").Value = 11.2
$ws.Range("AJ4").Value = 82458180
$ws.Range("D5").Value = 60291
$ws.Range("E5").Value = 7315
$ws.Range("F5").Value = 7315
$ws.Range("G5").Value = 6983
$ws.Range("H5").Value = 4895
$ws.Range("I5").Value = 2026
$ws.Range("J5").Value = 2869
$ws.Range("K5").Value = 73352
$ws.Range("L5").Value = 14033
$ws.Range("M5").Value = 59320
$ws.Range("N5").Value = 29973
$ws.Range("O5").Value = 29347
$ws.Range("P5").Value = 445
$ws.Range("Q5").Value = 5897
$ws.Range("R5").Value = -2435
$ws.Range("S5").Value = -873
$ws.Range("T5").Value = 8466
$ws.Range("U5").Value = -2568
$ws.Range("V5").Value = 2225
$ws.Range("W5").Value = 12.13
$ws.Range("X5").Value = 8.119999999999999
$ws.Range("Y5").Value = 6.95
$ws.Range("Z5").Value = 6.79
$ws.Range("AA5").Value = 23.66
$ws.Range("AB5").Value = 7014.91
$ws.Range("AC5").Value = 2279
$ws.Range("AD5").Value = 61.72
$ws.Range("AE5").Value = 35986
$ws.Range("AF5").Value = 3.91
$ws.Range("AG5").Value = 360
$ws.Range("AH5").Value = 0.26
$ws.Range("AI5").Value = 14.81
$ws.Range("AJ5").Value = 82458180
$ws.Range("D6").Value = 60782
$ws.Range("E6").Value = 5495
$ws.Range("F6").Value = 5495
$ws.Range("G6").Value = 5348
$ws.Range("H6").Value = 3763
$ws.Range("I6").Value = 1423
$ws.Range("K6").Value = 73874
$ws.Range("L6").Value = 11449
$ws.Range("M6").Value = 62425
$ws.Range("N6").Value = 31188
$ws.Range("P6").Value = 445
$ws.Range("Q6").Value = 7307
$ws.Range("R6").Value = -5404
$ws.Range("S6").Value = -1024
$ws.Range("T6").Value = 4369
$ws.Range("U6").Value = 2938
$ws.Range("V6").Value = 2141
$ws.Range("W6").Value = 9.039999999999999
$ws.Range("X6").Value = 6.19
$ws.Range("Y6").Value = 4.65
$ws.Range("Z6").Value = 5.11
$ws.Range("AA6").Value = 18.34
$ws.Range("AB6").Value = 7292.1
$ws.Range("AC6").Value = 1600
$ws.Range("AD6").Value = 45.32
$ws.Range("AE6").Value = 37446
$ws.Range("AF6").Value = 1.94
$ws.Range("AG6").Value = 310
$ws.Range("AH6").Value = 0.43
$ws.Range("AI6").Value = 18.17
$ws.Range("AJ6").Value = 82458180
$ws.Range("D7").Value = 62798
$ws.Range("E7").Value = 5164
$ws.Range("G7").Value = 5430
$ws.Range("H7").Value = 3956
$ws.Range("I7").Value = 1820
$ws.Range("K7").Value = 78638
$ws.Range("L7").Value = 13898
$ws.Range("M7").Value = 64741
$ws.Range("N7").Value = 32384
$ws.Range("P7").Value = 450
$ws.Range("Q7").Value = 7188
$ws.Range("R7").Value = -3710
$ws.Range("S7").Value = -904
$ws.Range("T7").Value = 3441
$ws.Range("U7").Value = 3034
$ws.Range("W7").Value = 8.220000000000001
$ws.Range("X7").Value = 6.3
$ws.Range("Y7").Value = 5.73
$ws.Range("Z7").Value = 5.19
$ws.Range("AA7").Value = 21.47
$ws.Range("AC7").Value = 2039
$ws.Range("AD7").Value = 35.42
$ws.Range("AE7").Value = 36093
$ws.Range("AF7").Value = 2
$ws.Range("AG7").Value = 360
$ws.Range("AH7").Value = 0.5
$ws.Range("AI7").Value = 16.31
$ws.Range("D8").Value = 67598
$ws.Range("E8").Value = 6695
$ws.Range("G8").Value = 6633
$ws.Range("H8").Value = 4848
$ws.Range("I8").Value = 2256
$ws.Range("K8").Value = 82323
$ws.Range("L8").Value = 13967
$ws.Range("M8").Value = 68356
$ws.Range("N8").Value = 34165
$ws.Range("P8").Value = 450
$ws.Range("Q8").Value = 7335
$ws.Range("R8").Value = -3632
$ws.Range("S8").Value = -573
$ws.Range("T8").Value = 3083
$ws.Range("U8").Value = 3852
$ws.Range("W8").Value = 9.9
$ws.Range("X8").Value = 7.17
$ws.Range("Y8").Value = 6.78
$ws.Range("Z8").Value = 6.02
$ws.Range("AA8").Value = 20.43
$ws.Range("AC8").Value = 2350
$ws.Range("AD8").Value = 30.73
$ws.Range("AE8").Value = 38078
$ws.Range("AF8").Value = 1.9
$ws.Range("AG8").Value = 418
$ws.Range("AH8").Value = 0.58
$ws.Range("AI8").Value = 15.26
$ws.Range("D9").Value = 72488
$ws.Range("E9").Value = 7910
$ws.Range("G9").Value = 8061
$ws.Range("H9").Value = 5880
$ws.Range("I9").Value = 2622
$ws.Range("K9").Value = 87707
$ws.Range("L9").Value = 13158
$ws.Range("M9").Value = 74548
$ws.Range("N9").Value = 36378
$ws.Range("P9").Value = 453
$ws.Range("Q9").Value = 8714
$ws.Range("R9").Value = -3124
$ws.Range("S9").Value = -721
$ws.Range("T9").Value = 2833
$ws.Range("U9").Value = 5409
$ws.Range("W9").Value = 10.91
$ws.Range("X9").Value = 8.109999999999999
$ws.Range("Y9").Value = 7.43
$ws.Range("Z9").Value = 6.92
$ws.Range("AA9").Value = 17.65
$ws.Range("AC9").Value = 2731
$ws.Range("AD9").Value = 26.43
$ws.Range("AE9").Value = 40545
$ws.Range("AF9").Value = 1.78
$ws.Range("AG9").Value = 480
$ws.Range("AH9").Value = 0.66
$ws.Range("AI9").Value = 15.1
